# Add two new line contingencies (line7, line8) to the lines_states sheet,
# inserted right after the existing "line6" row, and refresh the data for
# the "extr" rows that follow (now pushed down by two rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 8 (pushes former rows 8-15 down to 10-17).
$ws.Rows.Item(8).Resize(2).Insert()

# Apply the same formatting used by the other index cells (A column) onto
# the newly inserted rows: bold font, thin border, centered/top alignment.
foreach ($r in 8, 9) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop
    $cell.Borders.LineStyle = 1        # xlContinuous
    $cell.Borders.Weight = 2           # xlThin
}

# New rows: A=index, B=name, C=from_bus, D=to_bus, E=in_service
$newRows = @(
    @{ Row = 8;  A = 6; B = "line7"; C = 14; D = 11; E = $true },
    @{ Row = 9;  A = 7; B = "line8"; C = 16; D = 9;  E = $true }
)
foreach ($n in $newRows) {
    $ws.Cells.Item($n.Row, 1).Value = $n.A
    $ws.Cells.Item($n.Row, 2).Value = $n.B
    $ws.Cells.Item($n.Row, 3).Value = $n.C
    $ws.Cells.Item($n.Row, 4).Value = $n.D
    $ws.Cells.Item($n.Row, 5).Value = $n.E
}

# Updated data (index + from_bus/to_bus/in_service) for the rows that used
# to be 8-15 ("extr1".."extr8") and are now at rows 10-17.
$updates = @(
    @{ Row = 10; A = 8;  C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; C = 10; D = 11; E = $true  },
    @{ Row = 13; A = 11; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; C = 8;  D = 5;  E = $false }
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
